# Applies updated cryptocurrency price/volume figures (Price column D, Volume(1h) column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.271.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.984.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.99%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.37%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.977.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.15%  "
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("E10").Value = "  -5.93%  "
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("E13").Value = "  -5.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.13%  "
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.475.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.984.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.209.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.41%  "
$ws.Range("E19").Value = "  -4.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.661"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.65%  "
$ws.Range("E23").Value = "  -7.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.88%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -3.81%  "
$ws.Range("E29").Value = "  -4.51%  "
$ws.Range("E30").Value = "  -6.71%  "
$ws.Range("E31").Value = "  -8.71%  "
$ws.Range("E32").Value = "  -6.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0934"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.953"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "49.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0659"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.44%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  -6.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "384.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.632.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.16%  "
$ws.Range("E46").Value = "  -5.97%  "
$ws.Range("E47").Value = "  -5.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.77%  "
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.91%  "
